# Add a new "14-sep" data column (BT) to the right of the existing
# "13-sep" column (BS), carrying one value per data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header for the new column
$ws.Range("BT1").Value = "14-sep"

# Data values for rows 2-18
$ws.Range("BT2").Value = 0
$ws.Range("BT3").Value = 16.85517816169342
$ws.Range("BT4").Value = 16.958782744424354
$ws.Range("BT5").Value = 18.37033466968176
$ws.Range("BT6").Value = 0
$ws.Range("BT7").Value = 9.8420764294722094
$ws.Range("BT8").Value = 18.435433917128044
$ws.Range("BT9").Value = 12.160170348476889
$ws.Range("BT10").Value = 14.303058023483253
$ws.Range("BT11").Value = 14.76198637701328
$ws.Range("BT12").Value = 0
$ws.Range("BT13").Value = 9.8265049212000761
$ws.Range("BT14").Value = 0
$ws.Range("BT15").Value = 0
$ws.Range("BT16").Value = 7.4663191907695028
$ws.Range("BT17").Value = 0
$ws.Range("BT18").Value = 0

# Update selection to reflect the new active cell after entering data
$ws.Range("BV6").Select()
